# filter() method for featureGroupsSet and negate support for features method
#
# This progress-tracker sheet (fGroups) gets a new "done" (G) column mark for
# several methods that are now implemented for the Set variant, a couple of
# "almost as-is" (B) cells that flip from "X?" to a confirmed "X", a method
# (filter) that moves from "almost as-is" to "implement", groupInfo that
# moves from "implement" to "done", and two new remarks in column H about
# waiting for the autoID branch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cells marked "done" (column G) for methods now supported by the Set class ---
$doneCells = @(
    "G7", "G8", "G12", "G13", "G14", "G24",
    "G30", "G31", "G32", "G34", "G35", "G36", "G37", "G38",
    "G44", "G50", "G51", "G53"
)
foreach ($ref in $doneCells) {
    $ws.Range($ref).Value = "X"
}

# --- "almost as-is" (X?) confirmed as fully "as-is"/"implement" (X) ---
$confirmCells = @("B8", "B30", "B31", "B34", "B44", "B50", "B51")
foreach ($ref in $confirmCells) {
    $ws.Range($ref).Value = "X"
}

# --- featureTable: "as-is" (B13) moves to "implement" (C13) ---
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = "X"

# --- groupInfo: "implement" (C24) moves to "done" (G24) ---
$ws.Range("C24").ClearContents()

# --- Remarks: maybe wait for autoID branch ---
$ws.Range("H23").Value = "maybe wait for autoID branch"
$ws.Range("H48").Value = "maybe wait for autoID branch"

# Match the last-saved selection
$ws.Range("G15").Select() | Out-Null

Write-Output "edit applied"
